$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing "Local" F-column value; formulas F4/F5 recalc automatically
$ws.Range("F3").Value = 523287

# Leave the active cell on F10, matching the author's final selection
$ws.Range("F10").Select()
